# Update "Actual" hours on the Summary sheet for this week's status report
$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("F3").Value = 18
$summary.Range("F5").Value = 39
$summary.Range("F9").Value = 11
$summary.Range("F12").Value = 8
$summary.Range("F15").Value = 9

# Create this week's (W9) actuals on the Weekly sheet
$weekly = $wb.Worksheets.Item("Weekly")
$weekly.Range("K6").Value = 2
$weekly.Range("K8").Value = 3
$weekly.Range("K11").Value = 3
$weekly.Range("K14").Value = 3
$weekly.Range("K15").Value = 8
$weekly.Range("K16").Value = 2

# Leave the Weekly sheet active, selected on the last-entered cell
$weekly.Activate()
$weekly.Range("K14").Select()
